# Generate Report for Handoff
# Updates the localization-status workbook: the previously handed-back
# file pair is replaced by a newly-generated handoff (one .md source plus
# two .png dependency images), refreshing status/timestamps across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$mdName   = "24d8a804-661f-487c-bfec-a0035c828f78.md"
$png1Name = "ca2df28f-646e-406c-abe7-7ffdbf1a35fe.png"
$png2Name = "ca70ebeb-dbca-4eaa-9965-7a71af746acd.png"

$zhXlf = "24d8a804-661f-487c-bfec-a0035c828f78.36820d0340adc6f22c7a60e602be675e5327ed6b.zh-cn.xlf"
$deXlf = "24d8a804-661f-487c-bfec-a0035c828f78.36820d0340adc6f22c7a60e602be675e5327ed6b.de-de.xlf"
$png1Target = "a8133a05ba9521ad324460fa45be519bd1e0a178.png"
$png2Target = "a327a1e6d01fad86ddae9fdf3428bf251364cceb.png"

$status      = "Ready for handoff"
$overviewDate = "2016-03-23 07:21:57"
$zhHandoffDt  = "2016-03-23 07:21:49"
$deHandoffDt  = "2016-03-23 07:21:57"
$minDate      = "0001-01-01 00:00:00"
$dependencyFrom = "e2e\24d8a804-661f-487c-bfec-a0035c828f78.md"

# --------------------------------------------------------------------
# Sheet 1: Overview
# --------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/$mdName", "", "", $mdName)
$ws1.Range("B2").Value = $status
$ws1.Range("C2").Value = $status
$ws1.Range("D2").Value = $overviewDate

$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/$png1Name", "", "", $png1Name)
$ws1.Range("B3").Value = $status
$ws1.Range("C3").Value = $status
$ws1.Range("D3").Value = $overviewDate

$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/$png2Name", "", "", $png2Name)
$ws1.Range("B4").Value = $status
$ws1.Range("C4").Value = $status
$ws1.Range("D4").Value = $overviewDate

# --------------------------------------------------------------------
# Sheet 2: zh-cn
# --------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Hyperlinks.Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/$mdName", "", "", $mdName)
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = $status
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$zhXlf", "", "", $zhXlf)
$ws2.Range("E2").Value = $zhHandoffDt
$ws2.Range("F2").Clear()
$ws2.Range("G2").Clear()
$ws2.Range("H2").Value = $minDate
$ws2.Range("J2").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/$png1Name", "", "", $png1Name)
$ws2.Range("B3").Value = ".png"
$ws2.Range("C3").Value = $status
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$png1Target", "", "", $png1Target)
$ws2.Range("E3").Value = $zhHandoffDt
$ws2.Range("F3").Clear()
$ws2.Range("G3").Clear()
$ws2.Range("H3").Value = $minDate
$ws2.Range("J3").Value = "IsDependency"
$ws2.Range("K3").Value = $dependencyFrom

$ws2.Range("A4").Value = ""
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/$png2Name", "", "", $png2Name)
$ws2.Range("B4").Value = ".png"
$ws2.Range("C4").Value = $status
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$png2Target", "", "", $png2Target)
$ws2.Range("E4").Value = $zhHandoffDt
$ws2.Range("H4").Value = $minDate
$ws2.Range("J4").Value = "IsDependency"
$ws2.Range("K4").Value = $dependencyFrom

# --------------------------------------------------------------------
# Sheet 3: de-de
# --------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Hyperlinks.Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/$mdName", "", "", $mdName)
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = $status
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$deXlf", "", "", $deXlf)
$ws3.Range("E2").Value = $deHandoffDt
$ws3.Range("F2").Clear()
$ws3.Range("G2").Clear()
$ws3.Range("H2").Value = $minDate
$ws3.Range("J2").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/$png1Name", "", "", $png1Name)
$ws3.Range("B3").Value = ".png"
$ws3.Range("C3").Value = $status
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$png1Target", "", "", $png1Target)
$ws3.Range("E3").Value = $deHandoffDt
$ws3.Range("F3").Clear()
$ws3.Range("G3").Clear()
$ws3.Range("H3").Value = $minDate
$ws3.Range("J3").Value = "IsDependency"
$ws3.Range("K3").Value = $dependencyFrom

$ws3.Range("A4").Value = ""
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/$png2Name", "", "", $png2Name)
$ws3.Range("B4").Value = ".png"
$ws3.Range("C4").Value = $status
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$png2Target", "", "", $png2Target)
$ws3.Range("E4").Value = $deHandoffDt
$ws3.Range("H4").Value = $minDate
$ws3.Range("J4").Value = "IsDependency"
$ws3.Range("K4").Value = $dependencyFrom
